$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 10 values (group standings table)
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 2

# Update row 12 values
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = -4
$ws.Range("G12").Value = 1

# Update the active selection to match the saved view state
$ws.Range("H13").Select()
